$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.741.30'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '1.879.04'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.25%  '
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4711'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.96%  '
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.99'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08032'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.42%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = '1.895.75'
$ws.Range("E13").Value = '  +2.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.970'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.162'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06609'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Value = '27.700.37'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.507'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.296'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = '2.097.69'
$ws.Range("E26").Value = '  +1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.094'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.594'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9724'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09554'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.455'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.34%  '
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.309'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06120'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02266'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.235'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.165'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6002'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.249'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5686'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.400'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.933'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06823'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000315'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.71%  '
